$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $val) {
    $cell = $ws.Range($rangeAddr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue "D2" "288.68"
Set-TextValue "E2" "-6.98%"
Set-TextValue "D3" "40.18"
Set-TextValue "E3" "-2.37%"
Set-TextValue "D4" "5.032"
Set-TextValue "E4" "-3.34%"
Set-TextValue "D5" "0.07320"
Set-TextValue "E5" "-4.85%"
Set-TextValue "D6" "4.287"
Set-TextValue "E6" "-0.21%"
Set-TextValue "D7" "1.545"
Set-TextValue "E7" "-8.83%"
Set-TextValue "D8" "0.9110"
Set-TextValue "E8" "-3.52%"
Set-TextValue "D9" "0.1198"
Set-TextValue "E9" "-5.16%"
Set-TextValue "D10" "0.1737"
Set-TextValue "E10" "-5.28%"
Set-TextValue "D11" "0.08677"
Set-TextValue "E11" "-4.65%"
Set-TextValue "D12" "0.04163"
Set-TextValue "E12" "-1.75%"
Set-TextValue "D13" "0.1053"
Set-TextValue "E13" "0.12%"
Set-TextValue "D14" "0.001280"
Set-TextValue "E14" "-0.35%"
Set-TextValue "B15" "CoinExToken"
Set-TextValue "C15" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D15" "0.03838"
Set-TextValue "E15" "-4.54%"
Set-TextValue "B16" "TigerCash"
Set-TextValue "C16" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D16" "0.005795"
Set-TextValue "E16" "-1.58%"
Set-TextValue "B17" "LEO"
Set-TextValue "C17" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D17" "3.397"
Set-TextValue "E17" "1.39%"
Set-TextValue "B18" "BTSEToken"
Set-TextValue "C18" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D18" "2.397"
Set-TextValue "E18" "-1.16%"
Set-TextValue "B19" "BitpandaEcosystemToken"
Set-TextValue "C19" "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue "D19" "0.3284"
Set-TextValue "E19" "-1.08%"
Set-TextValue "B20" "MCDex"
Set-TextValue "C20" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D20" "7.588"
Set-TextValue "E20" "1.78%"
Set-TextValue "B21" "ProBitToken"
Set-TextValue "C21" "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue "D21" "0.1340"
Set-TextValue "E21" "-0.96%"
Set-TextValue "B22" "ZBToken"
Set-TextValue "C22" "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-TextValue "D22" "0.2886"
Set-TextValue "E22" "6.05%"
Set-TextValue "D23" "0.001269"
Set-TextValue "E23" "0.18%"
Set-TextValue "D24" "0.003672"
Set-TextValue "E24" "-13.28%"
Set-TextValue "D25" "0.0001283"
Set-TextValue "E25" "0.99%"
Set-TextValue "D26" "0.0003730"
Set-TextValue "E26" "-95.03%"
Set-TextValue "D38" "0.02327"
Set-TextValue "E38" "-8.19%"
Set-TextValue "D39" "0.04998"
Set-TextValue "E39" "-6.12%"
Set-TextValue "D40" "0.007684"
Set-TextValue "E40" "-1.98%"
Set-TextValue "E41" "154.29%"
Set-TextValue "D42" "0.1271"
Set-TextValue "E42" "-3.19%"
Set-TextValue "D43" "0.007386"
Set-TextValue "E43" "11.32%"
Set-TextValue "D44" "0.007000"
Set-TextValue "E44" "-5.70%"
Set-TextValue "D45" "0.3132"
Set-TextValue "E45" "1.35%"
Set-TextValue "D46" "0.00006539"
Set-TextValue "E46" "-3.64%"
Set-TextValue "D47" "0.00000000752"
Set-TextValue "E47" "0.10%"
Set-TextValue "E48" "12.29%"
Set-TextValue "D49" "0.004202"
Set-TextValue "E49" "35.38%"
Set-TextValue "D50" "0.00002105"
Set-TextValue "E50" "0.10%"
Set-TextValue "D51" "0.0002004"
Set-TextValue "E51" "0.10%"
